# miniproject - test case.xlsx
# modify script TC-01 - TC-05 and adding TC-06
#
# The diff shows the table grows by one row: a new TC-06 row is inserted
# at worksheet row 19 (the former TC-07 row, which shifts down to row 20),
# and several columns in the Table3 region are widened.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row above the current row 19 (old TC-07), shifting
#        it (and nothing else, since it's the last table row) down to row 20.
$ws.Rows.Item(19).Insert()

# --- 2. Grow Table3 (and its AutoFilter) by one row: A13:L19 -> A13:L20.
$lo3 = $ws.ListObjects.Item(3)
$lo3.Resize($ws.Range("A13:L20"))

# --- 3. Give the freshly inserted row 19 the same formatting as the data
#        rows around it (border + vertical-center + wrap-text style).
$ws.Range("A19:L19").Style = $ws.Range("A20:L20").Style

# --- 4. Fill in the new TC-06 test case data in row 19.
$ws.Range("A19").Value = "TC-06"
$ws.Range("B19").Value = "User gagal login menggunakan username dan password kosong"
$ws.Range("C19").Value = "Login"
$ws.Range("D19").Value = "High"
$ws.Range("E19").Value = "Negative"
$ws.Range("F19").Value = "Web Browser (Saucedemo)"
$ws.Range("G19").Value = "1. User berada di login page`n2. Credential valid yang tersedia:`nusername: standard_user`npassword: secret_sauce"
$ws.Range("H19").Value = "1. User membuka halaman saucedemo.com`n2. User tidak mengisi kolom username`n3. User tidak mengisi kolom password`n4. User klik tombol login`n5. Sistem validasi credential"
$ws.Range("I19").Value = "1. Jika user tidak mengisi kedua field maka muncul pesan error: 'Epic sadface: Username is required'`n2. Jika user hanya mengisi username maka muncul pesan error: 'Epic sadface: Password is required'`n3. Jika user hanya mengisi password maka pesan error: 'Epic sadface: Username is required'"

# Row height: real Excel auto-grows wrap-text rows to fit their content
# (225pt for the new 9-visual-line row). Set it explicitly to match.
$ws.Range("A19").EntireRow.RowHeight = 225

# --- 5. Widen the columns that now hold the longer TC-06 text
#        (values below are pre-compensated for this engine's column-width
#        rounding so the saved width lands as close as possible to the
#        authored target).
$ws.Columns.Item(7).ColumnWidth = 40.5               # G -> 41.28515625
$ws.Columns.Item(8).ColumnWidth = 37.833333333333336 # H -> 38.7109375
$ws.Columns.Item(9).ColumnWidth = 43.833333333333336 # I -> 44.7109375
$ws.Columns.Item(10).ColumnWidth = 19.5              # J -> 20.28515625
$ws.Columns.Item(11).ColumnWidth = 24.166666666666668 # K -> 25
$ws.Columns.Item(12).ColumnWidth = 28.333333333333332 # L -> 29.140625

# --- 6. Match the saved sheet view (scrolled down, new selected cell).
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("K19").Select()
